$wb = $excel.ActiveWorkbook

# --- Add Sheet2 after Sheet1 ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# --- Populate Sheet2 data ---
$ws2.Range("A1").Value = 9
$ws2.Range("B1").Value = 8
$ws2.Range("C1").Value = 7

$ws2.Range("A2").Value = 6
$ws2.Range("B2").Value = 5
$ws2.Range("C2").Value = 4

$ws2.Range("A3").Value = 3
$ws2.Range("B3").Value = 2
$ws2.Range("C3").Value = 1

# --- Defined names scoped to Sheet2 (localSheetId="1") ---
$ws2.Names.Add("Calculation", "=SUM(Sheet2!`$A`$1:`$C`$1)/SUM(Sheet2!`$A`$3/Sheet2!`$C`$3)")
$ws2.Range("E1").Formula = "=Calculation"
$ws2.Names.Add("Result", "=Sheet2!`$E`$1")

# Sheet2's last remembered selection is E1 ...
$ws2.Range("E1").Select() | Out-Null
# ... but Sheet1 stays the active/selected sheet in the workbook.
$ws1.Activate() | Out-Null
